$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same event rows and both
# need their F2/F3 "想去人数" (want-to-go count) values bumped.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 130
    $ws.Range("F3").Value = 32
}
